$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rename the "deals" sheet to "companies" ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "companies"

# --- Populate the companies sheet with header + data rows (entered column by
#     column, the way the source spreadsheet's data was pasted in) ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "industry"
$ws.Range("C1").Value = "employeeCount"
$ws.Range("D1").Value = "status"
$ws.Range("E1").Value = "category"

$ws.Range("A2").Value = "Company A"
$ws.Range("A3").Value = "Company B"
$ws.Range("A4").Value = "Company C"

$ws.Range("B2").Value = "IT"
$ws.Range("B3").Value = "IT"
$ws.Range("B4").Value = "IT"

$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 2000
$ws.Range("C4").Value = 3000

$ws.Range("D2").Value = "New"
$ws.Range("D3").Value = "Active"
$ws.Range("D4").Value = "Hot"

$ws.Range("E2").Value = "Client"
$ws.Range("E3").Value = "Client"
$ws.Range("E4").Value = "Partner"

# --- Header formatting: reuse the bold/yellow-fill look already used on the
#     "contacts" header row, copied over in one shot so no stray style gets
#     created ---
$ws1.Range("A1:D1").Copy() | Out-Null
$ws.Range("A1:E1").PasteSpecial(-4122) | Out-Null

# --- Widen the employeeCount column so the header text fits ---
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# --- Update the selection on the old active sheet (contacts) ---
$ws1.Range("A1:D4").Select() | Out-Null

# --- Make "companies" the active sheet with a fresh selection ---
$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null
